# Update 2p3. Added templates for formula student suspension, torque
# vectoring, four-wheel steering.
#
# The existing "Sedan_HambaLG" instance sheet is duplicated into a new
# "FSAE_Achilles" instance sheet (same BodyGeometry layout/styles), with
# the Instance label and rWheelCutout value updated for the new car.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet

# Duplicate the existing sheet right after itself.
$ws1.Copy($null, $ws1)

$ws2 = $wb.Worksheets.Item($ws1.Index + 1)
$ws2.Name = "FSAE_Achilles"

# Row 3 ("Instance") should read the new sheet's name.
$ws2.Range("H3").Value = "FSAE_Achilles"

# Row 6 ("rWheelCutout") differs for the FSAE_Achilles instance.
$ws2.Range("H6").Value = 0.25

# The newly added sheet becomes the active/selected tab.
$ws2.Activate()
